# Updated cryptos list - applies Price (D) and Volume(1h) (E) changes per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.075.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5036"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3886"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09223"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.68%  "

$ws.Range("E10").Value = "  -3.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.371"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.77%  "

$ws.Range("E13").Value = "  -2.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.897.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.281"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001106"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06638"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("E20").Value = "  -2.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.201"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.131.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.317"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.114.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.534"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.072"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1054"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.70%  "

$ws.Range("E33").Value = "  -3.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.489"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06585"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.336"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02399"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2198"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.214"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6411"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.952"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.36%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6026"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.301"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.687"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.96%  "

$ws.Range("E51").Value = "  -2.05%  "

